$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to Text format so a numeric-looking string (e.g. "310.81")
    # is stored as text (matching the workbook's inline-string cells) rather than
    # being auto-converted to a number by Excel's input parser. ClearFormats()
    # afterwards drops the now-unneeded explicit style so the cell keeps using the
    # default style index, leaving styles.xml / the cell's "s" attribute untouched.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "40.063.75"
Set-TextValue $ws.Range("E2") "  -2.95%  "
Set-TextValue $ws.Range("D3") "2.346.07"
Set-TextValue $ws.Range("E3") "  -3.72%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "310.81"
Set-TextValue $ws.Range("E5") "  -2.39%  "
Set-TextValue $ws.Range("D6") "84.98"
Set-TextValue $ws.Range("E6") "  -5.55%  "
Set-TextValue $ws.Range("D7") "0.526"
Set-TextValue $ws.Range("E7") "  -2.57%  "
Set-TextValue $ws.Range("E8") "  +0.03%  "
Set-TextValue $ws.Range("D9") "0.482"
Set-TextValue $ws.Range("E9") "  -3.88%  "
Set-TextValue $ws.Range("D10") "0.0806"
Set-TextValue $ws.Range("E10") "  -4.16%  "
Set-TextValue $ws.Range("D11") "30.05"
Set-TextValue $ws.Range("E11") "  -6.75%  "
Set-TextValue $ws.Range("E12") "  +0.55%  "
Set-TextValue $ws.Range("D13") "2.709.71"
Set-TextValue $ws.Range("E13") "  -3.34%  "
Set-TextValue $ws.Range("D14") "6.41"
Set-TextValue $ws.Range("E14") "  -5.05%  "
Set-TextValue $ws.Range("D15") "14.85"
Set-TextValue $ws.Range("E15") "  -5.66%  "
Set-TextValue $ws.Range("D16") "2.371.73"
Set-TextValue $ws.Range("E16") "  -1.97%  "
Set-TextValue $ws.Range("D17") "0.759"
Set-TextValue $ws.Range("E17") "  -2.53%  "
Set-TextValue $ws.Range("D18") "40.076.17"
Set-TextValue $ws.Range("E18") "  -2.66%  "
Set-TextValue $ws.Range("D19") "0.0₃0900"
Set-TextValue $ws.Range("E19") "  -3.26%  "
Set-TextValue $ws.Range("E20") "  -3.29%  "
Set-TextValue $ws.Range("D21") "68.20"
Set-TextValue $ws.Range("E21") "  -5.07%  "
Set-TextValue $ws.Range("D22") "10.62"
Set-TextValue $ws.Range("E22") "  -4.82%  "
Set-TextValue $ws.Range("D23") "234.95"
Set-TextValue $ws.Range("E23") "  -0.78%  "
Set-TextValue $ws.Range("E24") "  -5.45%  "
Set-TextValue $ws.Range("E25") "  -0.11%  "
Set-TextValue $ws.Range("E26") "  -3.42%  "
Set-TextValue $ws.Range("D27") "23.64"
Set-TextValue $ws.Range("E27") "  -2.57%  "
Set-TextValue $ws.Range("D28") "2.13"
Set-TextValue $ws.Range("E28") "  -4.57%  "
Set-TextValue $ws.Range("D29") "9.25"
Set-TextValue $ws.Range("E29") "  -4.08%  "
Set-TextValue $ws.Range("D30") "34.81"
Set-TextValue $ws.Range("E30") "  -0.38%  "
Set-TextValue $ws.Range("D31") "153.68"
Set-TextValue $ws.Range("E31") "  -1.46%  "
Set-TextValue $ws.Range("E32") "  +0.01%  "
Set-TextValue $ws.Range("D33") "5.10"
Set-TextValue $ws.Range("E33") "  -3.84%  "
Set-TextValue $ws.Range("D35") "0.0718"
Set-TextValue $ws.Range("E35") "  -4.24%  "
Set-TextValue $ws.Range("E36") "  -0.97%  "
Set-TextValue $ws.Range("D37") "2.79"
Set-TextValue $ws.Range("E37") "  -6.40%  "
Set-TextValue $ws.Range("D38") "0.0991"
Set-TextValue $ws.Range("E38") "  -2.31%  "
Set-TextValue $ws.Range("D39") "15.65"
Set-TextValue $ws.Range("E39") "  -7.28%  "
Set-TextValue $ws.Range("E40") "  -4.25%  "
Set-TextValue $ws.Range("E41") "  -2.16%  "
Set-TextValue $ws.Range("D42") "1.968.13"
Set-TextValue $ws.Range("E42") "  -1.77%  "
Set-TextValue $ws.Range("E43") "  +0.58%  "
Set-TextValue $ws.Range("D44") "0.0265"
Set-TextValue $ws.Range("E44") "  -4.22%  "
Set-TextValue $ws.Range("D45") "17.45"
Set-TextValue $ws.Range("E45") "  -6.40%  "
Set-TextValue $ws.Range("D46") "9.49"
Set-TextValue $ws.Range("E46") "  -0.43%  "
Set-TextValue $ws.Range("D47") "2.68"
Set-TextValue $ws.Range("E47") "  -8.42%  "
Set-TextValue $ws.Range("D48") "2.569.64"
Set-TextValue $ws.Range("E48") "  -3.46%  "
Set-TextValue $ws.Range("D49") "93.13"
Set-TextValue $ws.Range("E49") "  -2.23%  "
Set-TextValue $ws.Range("D50") "70.36"
Set-TextValue $ws.Range("E50") "  -4.83%  "
Set-TextValue $ws.Range("D51") "50.14"
Set-TextValue $ws.Range("E51") "  -4.04%  "
